# Append two new daily rows (2025-09-30) to the charging-volume log.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row60 = 60
$row61 = 61

# --- Row 60: 四方坪站充电量(kw) ---
$ws.Cells.Item($row60, 1).Value = 45930
$ws.Cells.Item($row60, 2).Value = "四方坪站充电量(kw)"

$row60Values = @(802.48400000000004, 1289.6679999999999, 322.73, 568.67400000000009, 319.57599999999996, 1086.7180000000001, 516.76400000000001, 318.90299999999996, 120.36499999999999, 185.059, 305.01, 224.76100000000002, 729.596, 1496.3490000000002, 576.20499999999993, 535.8599999999999, 423.15800000000002, 325.97999999999996, 139.9, 253.36999999999998, 224.85999999999996, 72.738, 53.84, 86.343000000000004)

for ($i = 0; $i -lt $row60Values.Length; $i++) {
    $ws.Cells.Item($row60, 3 + $i).Value = $row60Values[$i]
}

# --- Row 61: 高岭站充电量(kw) ---
$ws.Cells.Item($row61, 1).Value = 45930
$ws.Cells.Item($row61, 2).Value = "高岭站充电量(kw)"

$row61Values = @(471.54399999999998, 352.48, 254.34899999999999, 102.94499999999999, 35.75, 256.024, 385.82400000000007, 303.59100000000001, 456.80400000000003, 235.15699999999998, 72.616000000000014, 303.95300000000003, 288.54100000000005, 332.07600000000002, 526.17500000000007, 260.428, 117.21000000000001, 88.77600000000001, 176.86600000000001, 89.134, 131.29599999999999, 3.57, 130.91400000000002, 55.676000000000002)

for ($i = 0; $i -lt $row61Values.Length; $i++) {
    $ws.Cells.Item($row61, 3 + $i).Value = $row61Values[$i]
}

# Apply formats matching the existing data rows: date column uses the yyyy-mm-dd
# style already used throughout column A, numeric columns get a new "0.00_ " format.
$ws.Range("A60:A61").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("C60:Z61").NumberFormat = "0.00_ "

# Update the view to reflect scrolling down to the newly appended rows.
$ws.Application.ActiveWindow.ScrollRow = 49
$ws.Range("E74").Select()
